$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Add a new "2021" data column (column R) to the 5.5.1 table, one year
# after the existing last column ("2020", column Q).
#
# Copy the formatting of the last existing year column (Q) onto the new
# column (R) before writing the new values, so the appended header/value
# cells pick up the same look (font/borders/number format/alignment) as
# the rest of the table row, instead of landing with no style at all.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial($xlPasteFormats)
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial($xlPasteFormats)
$ws.Range("R5").Value = 20.5

# Move the active selection the way it appears after the edit in the
# source file (single cell S12 instead of the old N12:N13 block).
[void]$ws.Range("S12").Select()
